$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 37 (shifts rows 37:47 down to 38:48)
$ws.Rows.Item(37).Insert()

# Fill in the new row 37 with the "Fahrzeug - Details" bundle
$ws.Range("B37").Value = "menu-bar"
$ws.Range("C37").Value = "tixi_fahrzeug_details_page"
$ws.Range("D37").Value = "tixi_fahrzeug_page"
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = "/app/fahrzeug/details"
$ws.Range("G37").Value = "Details"
$ws.Range("H37").Value = "Fahrzeug - Details"
$ws.Range("I37").Value = "Vertrauliche Daten zur Fahrzeug"
$ws.Range("J37").Value = "ROLE_ADMIN"

# Selection used by the author after the edit
$ws.Range("A39").Select()
